# Applies the odds-value corrections described in the commit diff
# ("Atualizando o arquivo XLSX") to the single worksheet "Sheet1".
# Each statement updates one previously-existing numeric cell in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 21).Value = 1.95  # U2: 1.92 -> 1.95
$ws.Cells.Item(2, 22).Value = 1.8  # V2: 1.77 -> 1.8
# Row 5
$ws.Cells.Item(5, 22).Value = 1.63  # V5: 1.67 -> 1.63
# Row 6
$ws.Cells.Item(6, 7).Value = 1.85  # G6: 1.8 -> 1.85
$ws.Cells.Item(6, 8).Value = 3.1  # H6: 3.2 -> 3.1
$ws.Cells.Item(6, 9).Value = 4.75  # I6: 5 -> 4.75
$ws.Cells.Item(6, 10).Value = 2.6  # J6: 2.5 -> 2.6
$ws.Cells.Item(6, 14).Value = 7.5  # N6: 8 -> 7.5
$ws.Cells.Item(6, 22).Value = 1.63  # V6: 1.67 -> 1.63
$ws.Cells.Item(6, 24).Value = 8  # X6: 7.5 -> 8
$ws.Cells.Item(6, 34).Value = 10  # AH6: 11 -> 10
$ws.Cells.Item(6, 35).Value = 21  # AI6: 23 -> 21
$ws.Cells.Item(6, 40).Value = 3.75  # AN6: 3.6 -> 3.75
$ws.Cells.Item(6, 43).Value = 41  # AQ6: 34 -> 41
$ws.Cells.Item(6, 44).Value = 67  # AR6: 51 -> 67
$ws.Cells.Item(6, 51).Value = 34  # AY6: 41 -> 34
# Row 7
$ws.Cells.Item(7, 21).Value = 1.67  # U7: 1.7 -> 1.67
# Row 8
$ws.Cells.Item(8, 13).Value = 1.05  # M8: 1.03 -> 1.05
$ws.Cells.Item(8, 15).Value = 1.29  # O8: 1.25 -> 1.29
$ws.Cells.Item(8, 21).Value = 1.77  # U8: 1.8 -> 1.77
$ws.Cells.Item(8, 22).Value = 1.87  # V8: 1.91 -> 1.87
# Row 9
$ws.Cells.Item(9, 13).Value = 1.04  # M9: 1.03 -> 1.04
$ws.Cells.Item(9, 15).Value = 1.22  # O9: 1.19 -> 1.22
# Row 14
$ws.Cells.Item(14, 7).Value = 4  # G14: 3.6 -> 4
$ws.Cells.Item(14, 8).Value = 4.33  # H14: 4.2 -> 4.33
$ws.Cells.Item(14, 9).Value = 1.7  # I14: 1.8 -> 1.7
$ws.Cells.Item(14, 10).Value = 4  # J14: 3.75 -> 4
$ws.Cells.Item(14, 11).Value = 2.63  # K14: 2.6 -> 2.63
$ws.Cells.Item(14, 12).Value = 2.2  # L14: 2.3 -> 2.2
$ws.Cells.Item(14, 17).Value = 1.36  # Q14: 1.4 -> 1.36
$ws.Cells.Item(14, 18).Value = 3.1  # R14: 2.88 -> 3.1
$ws.Cells.Item(14, 21).Value = 1.44  # U14: 1.4 -> 1.44
$ws.Cells.Item(14, 22).Value = 2.63  # V14: 2.75 -> 2.63
$ws.Cells.Item(14, 23).Value = 21  # W14: 19 -> 21
$ws.Cells.Item(14, 24).Value = 29  # X14: 23 -> 29
$ws.Cells.Item(14, 25).Value = 15  # Y14: 13 -> 15
$ws.Cells.Item(14, 27).Value = 26  # AA14: 23 -> 26
$ws.Cells.Item(14, 28).Value = 26  # AB14: 23 -> 26
$ws.Cells.Item(14, 30).Value = 9.5  # AD14: 9 -> 9.5
$ws.Cells.Item(14, 35).Value = 12  # AI14: 13 -> 12
$ws.Cells.Item(14, 37).Value = 15  # AK14: 17 -> 15
$ws.Cells.Item(14, 38).Value = 12  # AL14: 13 -> 12
$ws.Cells.Item(14, 41).Value = 21  # AO14: 19 -> 21
$ws.Cells.Item(14, 42).Value = 21  # AP14: 19 -> 21
$ws.Cells.Item(14, 45).Value = 101  # AS14: 81 -> 101
$ws.Cells.Item(14, 49).Value = 4.33  # AW14: 4.5 -> 4.33
$ws.Cells.Item(14, 50).Value = 8.5  # AX14: 9 -> 8.5
$ws.Cells.Item(14, 51).Value = 13  # AY14: 15 -> 13
$ws.Cells.Item(14, 52).Value = 23  # AZ14: 26 -> 23
$ws.Cells.Item(14, 53).Value = 34  # BA14: 41 -> 34
$ws.Cells.Item(14, 55).Value = 251  # BC14: 201 -> 251
# Row 17
$ws.Cells.Item(17, 17).Value = 1.41  # Q17: 1.44 -> 1.41
# Row 19
$ws.Cells.Item(19, 17).Value = 1.67  # Q19: 1.7 -> 1.67
# Row 20
$ws.Cells.Item(20, 17).Value = 1.47  # Q20: 1.5 -> 1.47
# Row 21
$ws.Cells.Item(21, 15).Value = 1.08  # O21: 1.1 -> 1.08
$ws.Cells.Item(21, 16).Value = 8  # P21: 7 -> 8
# Row 23
$ws.Cells.Item(23, 22).Value = 1.73  # V23: 1.69 -> 1.73
# Row 24
$ws.Cells.Item(24, 22).Value = 1.73  # V24: 1.69 -> 1.73
# Row 25
$ws.Cells.Item(25, 21).Value = 1.8  # U25: 1.77 -> 1.8
$ws.Cells.Item(25, 22).Value = 1.95  # V25: 1.92 -> 1.95
# Row 26
$ws.Cells.Item(26, 21).Value = 1.62  # U26: 1.58 -> 1.62
# Row 29
$ws.Cells.Item(29, 8).Value = 3.65  # H29: 3.6 -> 3.65
$ws.Cells.Item(29, 9).Value = 1.85  # I29: 1.87 -> 1.85
$ws.Cells.Item(29, 10).Value = 4.05  # J29: 4.1 -> 4.05
$ws.Cells.Item(29, 11).Value = 2.22  # K29: 2.18 -> 2.22
$ws.Cells.Item(29, 12).Value = 2.4  # L29: 2.47 -> 2.4
$ws.Cells.Item(29, 14).Value = 8  # N29: 7.8 -> 8
$ws.Cells.Item(29, 15).Value = 1.24  # O29: 1.26 -> 1.24
$ws.Cells.Item(29, 16).Value = 3.7  # P29: 3.5 -> 3.7
$ws.Cells.Item(29, 17).Value = 1.72  # Q29: 1.78 -> 1.72
$ws.Cells.Item(29, 18).Value = 2.05  # R29: 1.95 -> 2.05
$ws.Cells.Item(29, 19).Value = 1.35  # S29: 1.38 -> 1.35
$ws.Cells.Item(29, 20).Value = 2.95  # T29: 2.82 -> 2.95
$ws.Cells.Item(29, 21).Value = 1.65  # U29: 1.7 -> 1.65
$ws.Cells.Item(29, 22).Value = 2.1  # V29: 2.05 -> 2.1
$ws.Cells.Item(29, 23).Value = 12.5  # W29: 12 -> 12.5
$ws.Cells.Item(29, 24).Value = 22  # X29: 21 -> 22
$ws.Cells.Item(29, 29).Value = 8  # AC29: 7.8 -> 8
$ws.Cells.Item(29, 30).Value = 7.1  # AD29: 6.9 -> 7.1
$ws.Cells.Item(29, 33).Value = 350  # AG29: 400 -> 350
$ws.Cells.Item(29, 34).Value = 8.25  # AH29: 7.8 -> 8.25
$ws.Cells.Item(29, 35).Value = 9.5  # AI29: 9.25 -> 9.5
$ws.Cells.Item(29, 38).Value = 14  # AL29: 14.5 -> 14
$ws.Cells.Item(29, 39).Value = 23  # AM29: 24 -> 23
$ws.Cells.Item(29, 40).Value = 5.7  # AN29: 5.6 -> 5.7
$ws.Cells.Item(29, 42).Value = 25  # AP29: 26 -> 25
$ws.Cells.Item(29, 43).Value = 100  # AQ29: 110 -> 100
$ws.Cells.Item(29, 44).Value = 120  # AR29: 150 -> 120
$ws.Cells.Item(29, 46).Value = 2.95  # AT29: 2.82 -> 2.95
$ws.Cells.Item(29, 47).Value = 7  # AU29: 7.1 -> 7
$ws.Cells.Item(29, 48).Value = 55  # AV29: 60 -> 55
$ws.Cells.Item(29, 49).Value = 3.85  # AW29: 3.8 -> 3.85
$ws.Cells.Item(29, 50).Value = 9.25  # AX29: 9.5 -> 9.25
$ws.Cells.Item(29, 51).Value = 17  # AY29: 18 -> 17
$ws.Cells.Item(29, 52).Value = 32  # AZ29: 35 -> 32
$ws.Cells.Item(29, 53).Value = 60  # BA29: 65 -> 60
$ws.Cells.Item(29, 54).Value = 200  # BB29: 250 -> 200
# Row 30
$ws.Cells.Item(30, 7).Value = 1.96  # G30: 2 -> 1.96
$ws.Cells.Item(30, 21).Value = 1.77  # U30: 1.8 -> 1.77
$ws.Cells.Item(30, 22).Value = 1.87  # V30: 1.91 -> 1.87
# Row 32
$ws.Cells.Item(32, 7).Value = 2.32  # G32: 2.38 -> 2.32
$ws.Cells.Item(32, 9).Value = 2.65  # I32: 2.7 -> 2.65
$ws.Cells.Item(32, 14).Value = 15  # N32: 17 -> 15
# Row 33
$ws.Cells.Item(33, 9).Value = 2.25  # I33: 2.3 -> 2.25
# Row 34
$ws.Cells.Item(34, 7).Value = 1.58  # G34: 1.62 -> 1.58
$ws.Cells.Item(34, 56).Value = 151  # BD34: 126 -> 151
# Row 36
$ws.Cells.Item(36, 21).Value = 1.73  # U36: 1.69 -> 1.73
# Row 37
$ws.Cells.Item(37, 13).Value = 1.04  # M37: 1.05 -> 1.04
$ws.Cells.Item(37, 14).Value = 9  # N37: 8.85 -> 9
# Row 38
$ws.Cells.Item(38, 22).Value = 1.73  # V38: 1.69 -> 1.73
# Row 39
$ws.Cells.Item(39, 17).Value = 1.95  # Q39: 1.93 -> 1.95
$ws.Cells.Item(39, 18).Value = 1.9  # R39: 1.93 -> 1.9
$ws.Cells.Item(39, 21).Value = 1.73  # U39: 1.69 -> 1.73
